$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text values that
# must remain stored as text (not auto-converted to number/percentage).
# Force the whole data range to Text format before writing the new values,
# then reset the style back to Normal so no stray number-format style sticks.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @{
    "D2" = '260.93'
    "E2" = '1.81%'
    "D3" = '27.22'
    "E3" = '1.40%'
    "D4" = '4.736'
    "E4" = '4.26%'
    "D5" = '0.06084'
    "E5" = '3.39%'
    "D6" = '6.666'
    "E6" = '0.91%'
    "D7" = '0.8456'
    "E7" = '-0.54%'
    "D8" = '0.9225'
    "E8" = '-0.62%'
    "D9" = '0.1406'
    "E9" = '2.24%'
    "D10" = '0.04981'
    "E10" = '11.41%'
    "D11" = '0.07091'
    "E11" = '1.06%'
    "D12" = '0.03133'
    "E12" = '1.59%'
    "D13" = '0.09078'
    "E13" = '-0.27%'
    "D14" = '0.001533'
    "E14" = '0.25%'
    "D15" = '0.0006075'
    "E15" = '-94.11%'
    "D16" = '0.006181'
    "E16" = '-0.58%'
    "D17" = '3.455'
    "E17" = '-0.72%'
    "D18" = '3.148'
    "E18" = '-0.61%'
    "D19" = '2.168'
    "E19" = '-1.65%'
    "E20" = '2.58%'
    "E21" = '0.86%'
    "D22" = '4.090'
    "E22" = '4.66%'
    "D23" = '0.04248'
    "E23" = '-0.58%'
    "D24" = '0.001221'
    "E24" = '0.04%'
    "E25" = '-9.00%'
    "E26" = '0.00%'
    "D27" = '0.0001574'
    "E27" = '3.36%'
    "D40" = '0.03871'
    "E40" = '1.82%'
    "D41" = '0.1116'
    "E41" = '1.66%'
    "D42" = '0.004128'
    "E42" = '-34.05%'
    "D43" = '0.01635'
    "E43" = '18.05%'
    "E44" = '0.36%'
    "D45" = '0.00005318'
    "E45" = '-0.64%'
    "E46" = '0.00%'
    "B47" = 'BOLO'
    "C47" = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
    "D47" = '0.1358'
    "E47" = '-46.24%'
    "B48" = 'CoinbaseStockToken'
    "C48" = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
    "D48" = '0.05448'
    "E48" = '1.29%'
    "E49" = '0.00%'
    "E50" = '0.00%'
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# Clear the temporary text number-format so untouched cells keep their
# original default (General) style, matching the source workbook.
$ws.Range("D2:E51").Style = "Normal"

